# Update vm_pu results: column B set to 1.02 (was 1.05) and columns
# C:F / I:N recomputed to the new power-flow solution values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# One tuple per data row (rows 2-25), values for columns B,C,D,E,F
$bfRows = @(
    (1.02,1.05334205180682,1.05842899860705,1.059965286444523,1.070491200670201),
    (1.02,1.054341776837865,1.059219268414817,1.060850626133397,1.071447858155445),
    (1.02,1.054989013788294,1.059730907244744,1.061424161477216,1.072067651268543),
    (1.02,1.055261194868426,1.059946066488843,1.061665432981164,1.07232839593013),
    (1.02,1.055306900098322,1.059982196527615,1.061705952717716,1.072372186842127),
    (1.02,1.054992650361183,1.059733781954027,1.061427384742351,1.072071134631307),
    (1.02,1.053679841114616,1.058696015061717,1.060264353544211,1.070814347075593),
    (1.02,1.051369197255987,1.056869540688981,1.058220057124983,1.068605686959257),
    (1.02,1.04983061563817,1.055653436443718,1.056860694054191,1.067137319551857),
    (1.02,1.049164838114699,1.055127229531381,1.056272918280567,1.066502479627845),
    (1.02,1.048917605065462,1.054931830097431,1.056054719004713,1.066266819008109),
    (1.02,1.048970634382596,1.05497374133265,1.056101517735165,1.066317362333539),
    (1.02,1.049144400396999,1.055111076579605,1.056254879267499,1.066482996838795),
    (1.02,1.049251472182511,1.055195700973488,1.056349387203773,1.066585069259749),
    (1.02,1.049874810379613,1.055688367035193,1.056899720525556,1.067179472371946),
    (1.02,1.050265931698646,1.05599750435419,1.057245154986289,1.067552586930282),
    (1.02,1.050494108419564,1.056177854963834,1.057446721792986,1.067770312136352),
    (1.02,1.050571917940098,1.056239355943439,1.057515464527732,1.067844566683224),
    (1.02,1.050223963706559,1.055964333110176,1.057208084780217,1.067512545559938),
    (1.02,1.04909322880074,1.055070633184031,1.056209714609799,1.066434217533732),
    (1.02,1.048382675111437,1.054509060261116,1.055582733882708,1.065757082038923),
    (1.02,1.048759316480666,1.054806728913085,1.055915038296684,1.066115963258347),
    (1.02,1.050242927105043,1.055979321656219,1.057224834961301,1.067530638241888),
    (1.02,1.051966231111543,1.057341460392908,1.058747944285346,1.06917596549898)
)

# One tuple per data row (rows 2-25), values for columns I,J,K,L,M,N
$inRows = @(
    (1.050116027385218,1.058359900749752,1.061161257660803,1.06269335357262,1.073190896744092,1.023208494420248),
    (1.05042028922828,1.059010226379492,1.061765618989672,1.063392846585844,1.07396357915616,1.023429061581402),
    (1.050616018858009,1.059430752297112,1.062156304000965,1.063845486854482,1.074463700772439,1.023571590393608),
    (1.050698028260436,1.059607474026287,1.062320456890707,1.064035781121018,1.074673985530206,1.023631462988442),
    (1.050711781840642,1.059637142409982,1.062348013518736,1.064067732592953,1.074709295210966,1.023641513122877),
    (1.050617115754114,1.059433113927149,1.062158497778817,1.063848029559964,1.074466510476702,1.023572390597357),
    (1.05021909154421,1.058579738501786,1.061365582023365,1.062929745543235,1.073451997909044,1.023283075859971),
    (1.049508949452723,1.05707388405131,1.059965509246173,1.061311814154867,1.071665444977578,1.02277180398448),
    (1.049029650412557,1.056068614122918,1.059030254549401,1.06023337890562,1.0704752364921,1.022429992694378),
    (1.048820721605477,1.055633005729537,1.058624844078699,1.059766458685395,1.069960070422136,1.022281760196308),
    (1.048742907834626,1.055471153975814,1.05847419132705,1.059593031910102,1.069768746020929,1.022226666409566),
    (1.048759608572667,1.05550587386082,1.05850650979513,1.059630232159408,1.069809784328118,1.022238485731001),
    (1.048814293742049,1.055619627976237,1.058612392402285,1.059752123002642,1.069944254845703,1.022277206811096),
    (1.048847959488634,1.055689709379728,1.058677621545733,1.059827225081559,1.070027110731655,1.022301059698568),
    (1.049043487095151,1.056097517333898,1.059057151117757,1.060264367968367,1.070509430697836,1.022439825662956),
    (1.049165764758347,1.056353239311802,1.059295102951718,1.060538589908976,1.070812031971853,1.022526809610824),
    (1.049236953190026,1.056502366706991,1.059433853824656,1.060698543638297,1.070988553613461,1.022577524053878),
    (1.049261203854556,1.056553210032442,1.059481157072466,1.060753084502106,1.071048746202481,1.022594812657548),
    (1.049152659385929,1.056325805955876,1.059269577353007,1.060509167997107,1.070779563702131,1.022517479311334),
    (1.04879819606982,1.0555861315277,1.058581214409133,1.059716228955257,1.069904655761279,1.022265805344157),
    (1.048574125898839,1.055120794936983,1.058148036128133,1.059217724440559,1.06935474744001,1.022107373473134),
    (1.048693023824584,1.055367504412512,1.058377707647017,1.059481986334261,1.069646246799588,1.022191379534404),
    (1.049158581553843,1.056338202004087,1.059281111403389,1.060522462478114,1.070794234655626,1.022521695339502),
    (1.04969357500593,1.057463427654091,1.060327795090404,1.06173005911109,1.072127170040976,1.022904151273013)
)

$bfValues = New-Object 'object[,]' $bfRows.Count, 5
for ($i = 0; $i -lt $bfRows.Count; $i++) {
    for ($j = 0; $j -lt 5; $j++) {
        $bfValues[$i, $j] = $bfRows[$i][$j]
    }
}

$inValues = New-Object 'object[,]' $inRows.Count, 6
for ($i = 0; $i -lt $inRows.Count; $i++) {
    for ($j = 0; $j -lt 6; $j++) {
        $inValues[$i, $j] = $inRows[$i][$j]
    }
}

$ws.Range("B2:F25").Value = $bfValues
$ws.Range("I2:N25").Value = $inValues
